$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (1) trial numbers
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON / meanEMG) values; C2 and E2 are cleared (no longer have data)
$ws.Range("B2").Value = 11.467718355161836
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 16.976148379153372
$ws.Range("E2").ClearContents()

# Update row 3 (STR / legmaxROM) values
$ws.Range("B3").Value = 10.83486683656362
$ws.Range("C3").Value = -3.105531684919832
$ws.Range("D3").Value = 18.13267575692705
$ws.Range("E3").Value = -0.39966137945635438

# Reflect the user's active selection after the edit
$ws.Range("B1:E3").Select()
